# Edit summary (per the target commit):
#  1. Slide 6's table changes its Table Style (tableStyleId) from the
#     custom "Table_0" style {69D65531-3413-4EF5-ADEB-A3D6D3EB3199} to the
#     built-in style {88D28B11-A518-4752-A44F-9426EABE568B}.
#  2. The deck's theme colour scheme is swapped from the custom "Integral"
#     palette to the stock Office palette (dk2/lt2/accent1-6/hlink/folHlink
#     all change; dk1/lt1 stay black/white). Fonts and format scheme are
#     untouched.

$p = $ppt.ActivePresentation

# --- 1. Table style on the table (slide 6, shape 2) ------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{88D28B11-A518-4752-A44F-9426EABE568B}")
    }
}

# --- 2. Theme colour scheme -> stock "Office" palette -----------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
$themeColors.Item(1).RGB  = 0x000000  # dk1
$themeColors.Item(2).RGB  = 0xFFFFFF  # lt1
$themeColors.Item(3).RGB  = 0x6A5444  # dk2      (44546A)
$themeColors.Item(4).RGB  = 0xE6E6E7  # lt2      (E7E6E6)
$themeColors.Item(5).RGB  = 0xD59B5B  # accent1  (5B9BD5)
$themeColors.Item(6).RGB  = 0x317DED  # accent2  (ED7D31)
$themeColors.Item(7).RGB  = 0xA5A5A5  # accent3  (A5A5A5)
$themeColors.Item(8).RGB  = 0x00C0FF  # accent4  (FFC000)
$themeColors.Item(9).RGB  = 0xC47244  # accent5  (4472C4)
$themeColors.Item(10).RGB = 0x47AD70  # accent6  (70AD47)
$themeColors.Item(11).RGB = 0xC16305  # hlink    (0563C1)
$themeColors.Item(12).RGB = 0x724F95  # folHlink (954F72)
